{"js": "// Append a SHA-256 verification footer block to the end of the document body.\nconst body = context.document.body;\n\nconst divider = \"\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\";\nconst lines = [\n  \"\",\n  divider,\n  \"\ud83d\udd10 SHA-256 Verification:\",\n  \"c387e9512105b8ea46fa32d895d23ef752adb3636d18c9b0f22bcdf6f7d5c4b8\",\n  \"Filed under: scroll-of-fire / 2_Witness_Scrolls/Book_of_Remembrance_and_Power.docx\",\n  divider,\n];\n\nfor (const line of lines) {\n  body.insertParagraph(line, \"End\");\n}\n\nawait context.sync();\n", "ps1": "# Append a SHA-256 verification footer block to the end of the document.\n$sel = $word.Selection\n$sel.EndKey(6)  # wdStory - jump to the very end of the main document story\n\n$divider = \"\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\u2500\"\n\n$sel.TypeParagraph()\n$sel.TypeParagraph()\n$sel.TypeText($divider)\n$sel.TypeParagraph()\n$sel.TypeText(\"\ud83d\udd10 SHA-256 Verification:\")\n$sel.TypeParagraph()\n$sel.TypeText(\"c387e9512105b8ea46fa32d895d23ef752adb3636d18c9b0f22bcdf6f7d5c4b8\")\n$sel.TypeParagraph()\n$sel.TypeText(\"Filed under: scroll-of-fire / 2_Witness_Scrolls/Book_of_Remembrance_and_Power.docx\")\n$sel.TypeParagraph()\n$sel.TypeText($divider)\n"}
